# "产出平衡" (output balance) sheet — dial in AI yield-bonus multipliers.
# Rows 18-19: columns C (瓶/science-ish), E (钱/gold-ish), F (粮/food-ish) get nudged.
# Rows 20-23: only column E changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C18").Value = 1.2
$ws.Range("E18").Value = 1.25
$ws.Range("F18").Value = 1.2

$ws.Range("C19").Value = 1.3
$ws.Range("E19").Value = 1.35
$ws.Range("F19").Value = 1.3

$ws.Range("E20").Value = 1.45
$ws.Range("E21").Value = 1.55
$ws.Range("E22").Value = 1.65
$ws.Range("E23").Value = 1.75

# Match the author's final view state: scrolled down a bit, selection on L27.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("L27").Select()
